# Aufgabenliste geupdated, MethCopy eingefügt (fails)
#
# - Row 20 (static-Attribute...): marker column B changes from "r" to "t"
# - Row 22 (Shallow-Copy & Deep-Copy): marker column B gets "t"
# - Row 19 (Konstruktoren...): filename column E gets "MethConstr"
# - Row 21 (Objektidentität): filename column E gets "MethEquals"
# - Selection / active cell moves to I20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B20").Value = "t"
$ws.Range("B22").Value = "t"
$ws.Range("E19").Value = "MethConstr"
$ws.Range("E21").Value = "MethEquals"

$ws.Range("I20").Select() | Out-Null
